# Updates the crypto price/volume figures to the latest scrape, and
# refreshes row 51 (replaced coin entry) per the upstream GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text such as "65.704.00"
# or "  +0.91%  " (European-style thousands separators, literal percent
# strings with padding). Force text format first so Excel's COM layer
# does not reinterpret them as numbers/dates when the .Value is assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '65.704.00'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.582.39'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '602.60'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").Value = '137.61'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '3.580.88'
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '7.20'
$ws.Range("E11").Value = '  +6.02%  '
$ws.Range("D12").Value = '0.391'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '4.197.27'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").Value = '28.12'
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '3.584.64'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D18").Value = '65.808.20'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '9.99'
$ws.Range("E19").Value = '  -2.60%  '
$ws.Range("D20").Value = '14.60'
$ws.Range("E20").Value = '  +2.34%  '
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Value = '395.95'
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("D23").Value = '0.590'
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = '3.729.34'
$ws.Range("E24").Value = '  +1.72%  '
$ws.Range("D25").Value = '74.30'
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +3.47%  '
$ws.Range("D28").Value = '8.10'
$ws.Range("E28").Value = '  +5.66%  '
$ws.Range("E29").Value = '  +31.39%  '
$ws.Range("E30").Value = '  +4.65%  '
$ws.Range("D31").Value = '8.58'
$ws.Range("E31").Value = '  +5.28%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '3.589.43'
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").Value = '24.49'
$ws.Range("E34").Value = '  +3.14%  '
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("D37").Value = '5.40'
$ws.Range("E37").Value = '  +9.00%  '
$ws.Range("E38").Value = '  +5.26%  '
$ws.Range("D39").Value = '7.07'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").Value = '168.99'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '0.0837'
$ws.Range("E41").Value = '  +5.00%  '
$ws.Range("D42").Value = '0.840'
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").Value = '27.10'
$ws.Range("E43").Value = '  +4.92%  '
$ws.Range("D44").Value = '1.27'
$ws.Range("E44").Value = '  +8.54%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  +2.23%  '
$ws.Range("D49").Value = '7.01'
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("D50").Value = '2.455.21'
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '0.904'
$ws.Range("E51").Value = '  +10.38%  '
